$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "8.4.2"
$ws.Range("C11").Value = "2020-03-08"
$ws.Range("C14").Value = "7.5"
$ws.Range("C15").Value = "20200312-675bb1f"
$ws.Range("C16").Value = "74.0"
$ws.Range("C23").Value = "44.0.0"
$ws.Range("C29").Value = "1.43"
$ws.Range("C30").Value = "5.17.2"
